$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relative Abundances")

# Fix capitalization of "C1 oxidation (Sum)" -> "C1 Oxidation (Sum)" for rows 9-11
$ws.Range("B9").Value = "C1 Oxidation (Sum)"
$ws.Range("B10").Value = "C1 Oxidation (Sum)"
$ws.Range("B11").Value = "C1 Oxidation (Sum)"

# Rows 12-13 remain "Iron Reduction" (unchanged value, underlying shared-string order changed only)
$ws.Range("B12").Value = "Iron Reduction"
$ws.Range("B13").Value = "Iron Reduction"

# Update the active selection from N11 to B11
$ws.Range("B11").Select()

$wb.Save()
